# "Correcion validacion fecha inventario"
# Fix: correct a miskeyed distributor code (row 5) and add several new
# distributors to the "Base Clientes carga manual" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Correct the distributor number in row 5 (was a typo/duplicate id) ---
$ws.Range("A5").Value = 500226

# --- Append the new distributor rows to the table ---
$newDistributors = @(
    @(62000076, "DIABONOS S.A."),
    @(10236216, "ARIANNA GARCIA"),
    @(10234501, "Agroquimicos Libra"),
    @(10220649, "ASESORIA INTEGRAL LUMINARIAS"),
    @(10234501, "Agroquimicos Libra")
)

foreach ($distributor in $newDistributors) {
    $row = $lo.ListRows.Add()
    $row.Range.Cells.Item(1, 1).Value = $distributor[0]
    $row.Range.Cells.Item(1, 2).Value = $distributor[1]
}

# --- Move the active selection, matching the author's final cursor spot ---
$ws.Range("C6").Select() | Out-Null
